$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the columns that are no longer part of the report.
# Delete from right-to-left so earlier column letters stay valid while iterating.
# Columns removed (original letters): V (HC Category), U (Phan loai),
# R (Phut nghi khong luong), P (Phut tang ca 150%), O (Phut tang ca 100%),
# N (Phut nghi phep), J (So phut ca), G (Cap bac)
$ws.Columns("V").Delete()
$ws.Columns("U").Delete()
$ws.Columns("R").Delete()
$ws.Columns("P").Delete()
$ws.Columns("O").Delete()
$ws.Columns("N").Delete()
$ws.Columns("J").Delete()
$ws.Columns("G").Delete()

# The old "Phut tang ca dem" column (now column L) is renamed to "Phut tang ca 200%"
$ws.Range("L3").Value = "Phút tăng ca 200%"

# Restore the selection that Excel recorded after the edit
$ws.Range("L7").Select()
